# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values per game row (column G), replacing old Strike# values
$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 4
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
